# Updated cryptos list on Tue May 16 07:55:56 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''27.344.79'
$ws.Range("E2").Value = '  -1.64%  '

# Row 3
$ws.Range("D3").Value = '''1.829.80'
$ws.Range("E3").Value = '  -1.37%  '

# Row 4
$ws.Range("E4").Value = '  -1.05%  '

# Row 5
$ws.Range("D5").Value = '''314.93'
$ws.Range("E5").Value = '  -1.84%  '

# Row 6
$ws.Range("E6").Value = '  -1.02%  '

# Row 7
$ws.Range("D7").Value = '''0.4279'
$ws.Range("E7").Value = '  -0.82%  '

# Row 8
$ws.Range("D8").Value = '''0.3701'
$ws.Range("E8").Value = '  -2.41%  '

# Row 9
$ws.Range("D9").Value = '''0.07259'
$ws.Range("E9").Value = '  -2.01%  '

# Row 10
$ws.Range("D10").Value = '''0.8673'
$ws.Range("E10").Value = '  -2.03%  '

# Row 11
$ws.Range("D11").Value = '''21.13'
$ws.Range("E11").Value = '  -2.61%  '

# Row 12
$ws.Range("D12").Value = '''1.828.37'
$ws.Range("E12").Value = '  -1.61%  '

# Row 13
$ws.Range("D13").Value = '''6.722'

# Row 14
$ws.Range("D14").Value = '''0.07114'
$ws.Range("E14").Value = '  +0.07%  '

# Row 15
$ws.Range("D15").Value = '''5.324'
$ws.Range("E15").Value = '  -3.04%  '

# Row 16
$ws.Range("D16").Value = '''88.99'
$ws.Range("E16").Value = '  +0.66%  '

# Row 17
$ws.Range("E17").Value = '  -1.20%  '

# Row 18
$ws.Range("D18").Value = '''0.000008875'
$ws.Range("E18").Value = '  -1.77%  '

# Row 19
$ws.Range("D19").Value = '''1.004'
$ws.Range("E19").Value = '  -0.96%  '

# Row 20
$ws.Range("D20").Value = '''15.11'
$ws.Range("E20").Value = '  -2.71%  '

# Row 21
$ws.Range("D21").Value = '''27.313.76'

# Row 22
$ws.Range("D22").Value = '''5.148'
$ws.Range("E22").Value = '  -2.52%  '

# Row 23
$ws.Range("D23").Value = '''10.89'
$ws.Range("E23").Value = '  -2.79%  '

# Row 24
$ws.Range("D24").Value = '''2.047.81'
$ws.Range("E24").Value = '  -2.13%  '

# Row 25
$ws.Range("D25").Value = '''2.008'
$ws.Range("E25").Value = '  -1.08%  '

# Row 26
$ws.Range("D26").Value = '''152.82'
$ws.Range("E26").Value = '  -2.34%  '

# Row 27
$ws.Range("D27").Value = '''2.198'
$ws.Range("E27").Value = '  +8.06%  '

# Row 28
$ws.Range("E28").Value = '  -1.17%  '

# Row 29
$ws.Range("D29").Value = '''5.260'
$ws.Range("E29").Value = '  -2.99%  '

# Row 30
$ws.Range("D30").Value = '''116.85'
$ws.Range("E30").Value = '  -3.85%  '

# Row 31
$ws.Range("D31").Value = '''0.08903'
$ws.Range("E31").Value = '  -0.78%  '

# Row 32
$ws.Range("E32").Value = '  -2.73%  '

# Row 33
$ws.Range("D33").Value = '''0.7609'
$ws.Range("E33").Value = '  -2.16%  '

# Row 34
$ws.Range("D34").Value = '''4.473'
$ws.Range("E34").Value = '  -2.32%  '

# Row 35
$ws.Range("D35").Value = '''2.818'
$ws.Range("E35").Value = '  -4.01%  '

# Row 36
$ws.Range("D36").Value = '''1.005'
$ws.Range("E36").Value = '  -1.05%  '

# Row 37
$ws.Range("E37").Value = '  -2.83%  '

# Row 38
$ws.Range("D38").Value = '''0.01981'
$ws.Range("E38").Value = '  +0.60%  '

# Row 39
$ws.Range("D39").Value = '''0.05290'
$ws.Range("E39").Value = '  -0.57%  '

# Row 40
$ws.Range("D40").Value = '''7.198'
$ws.Range("E40").Value = '  +2.77%  '

# Row 41
$ws.Range("D41").Value = '''2.887'
$ws.Range("E41").Value = '  +0.41%  '

# Row 42
$ws.Range("D42").Value = '''0.1702'
$ws.Range("E42").Value = '  +1.07%  '

# Row 43
$ws.Range("D43").Value = '''0.5088'
$ws.Range("E43").Value = '  -2.10%  '

# Row 44
$ws.Range("D44").Value = '''8.709'
$ws.Range("E44").Value = '  -0.91%  '

# Row 45
$ws.Range("D45").Value = '''10.69'
$ws.Range("E45").Value = '  -0.94%  '

# Row 46
$ws.Range("D46").Value = '''107.89'
$ws.Range("E46").Value = '  -2.54%  '

# Row 47
$ws.Range("D47").Value = '''0.4780'
$ws.Range("E47").Value = '  +0.74%  '

# Row 48
$ws.Range("D48").Value = '''1.005'
$ws.Range("E48").Value = '  -1.05%  '

# Row 49
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.06384'
$ws.Range("E49").Value = '  -2.19%  '

# Row 50
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '''1.671'
$ws.Range("E50").Value = '  -2.39%  '

# Row 51
$ws.Range("D51").Value = '''1.848'
$ws.Range("E51").Value = '  -1.76%  '
